$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells we touch to remain plain text so that
# values such as "0.9990", "13.08", "30.279.99" keep their exact
# formatting instead of being auto-coerced to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.290.34"
$ws.Range("E2").Value = "  +5.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.917.42"
$ws.Range("E3").Value = "  +5.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.12"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5169"
$ws.Range("E7").Value = "  +4.05%  "
$ws.Range("E8").Value = "  +6.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2978"
$ws.Range("E9").Value = "  +6.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06830"
$ws.Range("E10").Value = "  +6.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.918.28"
$ws.Range("E11").Value = "  +6.08%  "
$ws.Range("E12").Value = "  +4.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07331"
$ws.Range("E13").Value = "  +3.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6900"
$ws.Range("E14").Value = "  +6.21%  "
$ws.Range("E15").Value = "  +7.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.921"
$ws.Range("E16").Value = "  +4.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.277.09"
$ws.Range("E17").Value = "  +5.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007975"
$ws.Range("E18").Value = "  +8.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.08"
$ws.Range("E20").Value = "  +6.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.164.73"
$ws.Range("E21").Value = "  +6.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9989"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.863"
$ws.Range("E23").Value = "  +5.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.776"
$ws.Range("E24").Value = "  +8.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.211"
$ws.Range("E25").Value = "  +3.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.93"
$ws.Range("E26").Value = "  +2.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "139.80"
$ws.Range("E27").Value = "  +25.09%  "
$ws.Range("E28").Value = "  +7.82%  "
$ws.Range("E29").Value = "  +6.87%  "
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.289"
$ws.Range("E31").Value = "  +2.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08860"
$ws.Range("E32").Value = "  +5.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.031"
$ws.Range("E33").Value = "  +4.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05139"
$ws.Range("E34").Value = "  +3.06%  "
$ws.Range("E35").Value = "  +6.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7225"
$ws.Range("E36").Value = "  +6.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.681"
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9761"
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01704"
$ws.Range("E41").Value = "  +6.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.238"
$ws.Range("E42").Value = "  +5.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4340"
$ws.Range("E43").Value = "  +5.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.99"
$ws.Range("E44").Value = "  +4.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9990"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.703"
$ws.Range("E46").Value = "  +6.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1278"
$ws.Range("E47").Value = "  +4.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05732"
$ws.Range("E48").Value = "  +4.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3851"
$ws.Range("E51").Value = "  +6.73%  "

# Row 38/39 swap: RenderToken moves up to row 38, MXToken moves down to row 39
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.331"
$ws.Range("E38").Value = "  +8.34%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.851"
$ws.Range("E39").Value = "  +7.06%  "

# Row 49/50 swap: Elrond moves up to row 49, EnergySwap moves down to row 50
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.41"
$ws.Range("E49").Value = "  +6.27%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.545"
$ws.Range("E50").Value = "  +4.48%  "
